# Add a new date column "08-10-2020" (column W) to the COVID19 deceased
# cases timeseries sheet, with the corresponding data values for each
# state/UT (rows 2-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell W1 -------------------------------------------------
# Force the value to be treated as text (not auto-converted to a date
# serial number) by using a leading apostrophe, then mirror the
# formatting that is used by the other header cells (bold, centered,
# thin box border) so it matches the style used by column V ("07-10-2020").
$ws.Cells.Item(1, 23).Value = "'08-10-2020"
$ws.Cells.Item(1, 23).Font.Bold = $true
$ws.Cells.Item(1, 23).HorizontalAlignment = -4108   # xlCenter
$ws.Cells.Item(1, 23).VerticalAlignment = -4160     # xlTop
$ws.Cells.Item(1, 23).Borders.LineStyle = 1         # xlContinuous (thin box)

# --- Data values for column W, rows 2 through 36 --------------------
$values = @{
    2  = 54
    3  = 6086
    4  = 21
    5  = 785
    6  = 927
    7  = 182
    8  = 1134
    9  = 2
    10 = 5616
    11 = 477
    12 = 3531
    13 = 1528
    14 = 231
    15 = 1282
    16 = 767
    17 = 9574
    18 = 906
    19 = 63
    20 = 2518
    21 = 39072
    22 = 80
    23 = 60
    24 = 0
    25 = 17
    26 = 958
    27 = 551
    28 = 3712
    29 = 1590
    30 = 49
    31 = 9984
    32 = 1201
    33 = 304
    34 = 688
    35 = 6200
    36 = 5376
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 23).Value = $values[$row]
}

$wb.Save()
